# Update the Bilibili-view-count style numbers (column F) on the "展览"
# and "全部类型" worksheets to reflect the refreshed scrape output.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 1101
$ws1.Range("F6").Value  = 55
$ws1.Range("F7").Value  = 59
$ws1.Range("F8").Value  = 11317
$ws1.Range("F9").Value  = 4305
$ws1.Range("F13").Value = 2516
$ws1.Range("F18").Value = 495
$ws1.Range("F19").Value = 11261
$ws1.Range("F20").Value = 11123

# --- Sheet "全部类型" (all types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 1101
$ws4.Range("F6").Value  = 55
$ws4.Range("F7").Value  = 59
$ws4.Range("F8").Value  = 11317
$ws4.Range("F9").Value  = 4305
$ws4.Range("F11").Value = 27
$ws4.Range("F13").Value = 2516
$ws4.Range("F18").Value = 495
$ws4.Range("F19").Value = 11261
$ws4.Range("F20").Value = 11123
